$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$PasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Build three "template" cells, far away from the real data, that carry the
# --- exact new cell formats needed for the new 2021 / 2022 columns. Building
# --- them once (off to the side) and then stamping them via Copy/PasteSpecial
# --- keeps the resulting style table minimal (no throw-away intermediate xfs).

$tmplBold       = $ws.Range("BZ1")   # used for row 5  (N5/O5)
$tmplPlain      = $ws.Range("BZ2")   # used for rows 6-13 (N6:O13)
$tmplPlainBord  = $ws.Range("BZ3")   # used for row 14 (N14/O14)

$tmplBold.NumberFormat = "0.0"
$tmplBold.Font.Name = "Times New Roman"
$tmplBold.Font.Size = 10
$tmplBold.Font.Bold = $true

$tmplPlain.NumberFormat = "0.0"
$tmplPlain.Font.Name = "Times New Roman"
$tmplPlain.Font.Size = 10
$tmplPlain.Font.Bold = $false

# Seed the bottom-border template from an existing cell that already carries
# the medium bottom border (M3), so the shared border table is reused as-is
# instead of minting new (intermediate) border definitions. Only the number
# format still needs to change; the font it inherits (Times New Roman 10,
# not bold) is already what we want.
$ws.Range("M3").Copy()
$tmplPlainBord.PasteSpecial($PasteFormats)
$excel.CutCopyMode = $false
$tmplPlainBord.NumberFormat = "0.0"

# --- Row 3: trailing border-only cells (same style as the existing K3:M3). ---
$ws.Range("M3").Copy()
$ws.Range("N3:O3").PasteSpecial($PasteFormats)

# --- Row 4: header years 2021 / 2022, same style as the existing year cells. ---
$ws.Range("M4").Copy()
$ws.Range("N4:O4").PasteSpecial($PasteFormats)
$ws.Range("N4").Value = 2021
$ws.Range("O4").Value = 2022

# --- Row 5 (bold row) ---
$tmplBold.Copy()
$ws.Range("N5:O5").PasteSpecial($PasteFormats)
$ws.Range("N5").Value = 40.007977647471066
$ws.Range("O5").Value = 42.620582506455563

# --- Rows 6-13 (plain rows) ---
$tmplPlain.Copy()
$ws.Range("N6:O13").PasteSpecial($PasteFormats)

$ws.Range("N6").Value = 5.7072514621689896
$ws.Range("O6").Value = 8.1443914479075037

$ws.Range("N7").Value = 8.9893229854028949
$ws.Range("O7").Value = 10.715961386284755

$ws.Range("N8").Value = 66.307512472824584
$ws.Range("O8").Value = 81.977461999426666

$ws.Range("N9").Value = 23.475213049310256
$ws.Range("O9").Value = 29.828871240443185

$ws.Range("N10").Value = 9.8045372040896162
$ws.Range("O10").Value = 9.7218425128664112

$ws.Range("N11").Value = 9.3737779268960448
$ws.Range("O11").Value = 8.6167819403064012

$ws.Range("N12").Value = 70.457032471318783
$ws.Range("O12").Value = 69.915337594090886

$ws.Range("N13").Value = 98.411252120183207
$ws.Range("O13").Value = 99.08571752721997

# --- Row 14 (totals row, bottom border) ---
$tmplPlainBord.Copy()
$ws.Range("N14:O14").PasteSpecial($PasteFormats)
$ws.Range("N14").Value = 63.900563564170795
$ws.Range("O14").Value = 64.805252627098838

$excel.CutCopyMode = $false

# --- Remove the scratch template cells; they must leave no trace. ---
$ws.Range("BZ1:BZ3").Clear()

# --- Restore the selection to match the authored workbook. ---
$null = $ws.Range("P8").Select()
